$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 height shrinks slightly (layout adjustment that accompanies the edit)
$ws.Rows(3).RowHeight = 13.5

# New "2021" column (J) of data: copy the formatting from the adjacent
# 2020 column (I) for each data row, then set the new year's value.
$newValues = @{
    4  = 2021
    5  = 24.4
    6  = 45.7
    7  = 38
    8  = 51.3
    9  = 51.5
    10 = 13
    11 = 36.4
    12 = 27
    13 = 2.7
    14 = 40.4
}

foreach ($row in 4..14) {
    $src = $ws.Cells.Item($row, 9)
    $dst = $ws.Cells.Item($row, 10)
    $src.Copy($dst)
    $dst.Value = $newValues[$row]
}

# Selection ends up on K18 after the edits
$ws.Range("K18").Select()
